# Refresh cryptocurrency price/volume snapshot (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value. A leading apostrophe is used for
# numeric-looking 'Price' values so Excel stores them as text (preserving
# trailing zeros / exact formatting) instead of coercing them to numbers.
$updates = @(
    @{ Cell = 'D2'; Value = '35.282.17' }
    @{ Cell = 'E2'; Value = '  -0.91%  ' }
    @{ Cell = 'D3'; Value = '1.898.36' }
    @{ Cell = 'E3'; Value = '  -0.47%  ' }
    @{ Cell = 'E4'; Value = '  -0.18%  ' }
    @{ Cell = 'D5'; Value = '''0.691' }
    @{ Cell = 'E5'; Value = '  +8.93%  ' }
    @{ Cell = 'D6'; Value = '''245.30' }
    @{ Cell = 'E6'; Value = '  -0.07%  ' }
    @{ Cell = 'E7'; Value = '  -0.10%  ' }
    @{ Cell = 'D8'; Value = '''40.63' }
    @{ Cell = 'E8'; Value = '  -4.74%  ' }
    @{ Cell = 'E9'; Value = '  +2.17%  ' }
    @{ Cell = 'D10'; Value = '''53.16' }
    @{ Cell = 'E10'; Value = '  +11.19%  ' }
    @{ Cell = 'D11'; Value = '''0.0721' }
    @{ Cell = 'E11'; Value = '  +1.83%  ' }
    @{ Cell = 'E12'; Value = '  -0.01%  ' }
    @{ Cell = 'D13'; Value = '2.173.44' }
    @{ Cell = 'E13'; Value = '  -0.34%  ' }
    @{ Cell = 'D14'; Value = '''12.59' }
    @{ Cell = 'E14'; Value = '  +0.67%  ' }
    @{ Cell = 'D15'; Value = '''0.706' }
    @{ Cell = 'E15'; Value = '  +1.32%  ' }
    @{ Cell = 'D16'; Value = '1.897.57' }
    @{ Cell = 'E16'; Value = '  -1.38%  ' }
    @{ Cell = 'D17'; Value = '''4.80' }
    @{ Cell = 'E17'; Value = '  -0.74%  ' }
    @{ Cell = 'D18'; Value = '35.263.22' }
    @{ Cell = 'E18'; Value = '  -0.93%  ' }
    @{ Cell = 'D19'; Value = '''72.11' }
    @{ Cell = 'E19'; Value = '  +0.03%  ' }
    @{ Cell = 'D20'; Value = '0.0₃0817' }
    @{ Cell = 'E20'; Value = '  +0.82%  ' }
    @{ Cell = 'D21'; Value = '''240.69' }
    @{ Cell = 'E21'; Value = '  -1.60%  ' }
    @{ Cell = 'D22'; Value = '''12.59' }
    @{ Cell = 'E22'; Value = '  +1.07%  ' }
    @{ Cell = 'E24'; Value = '  -0.17%  ' }
    @{ Cell = 'E25'; Value = '  +1.22%  ' }
    @{ Cell = 'D26'; Value = '''2.29' }
    @{ Cell = 'E26'; Value = '  +7.54%  ' }
    @{ Cell = 'D27'; Value = '''167.66' }
    @{ Cell = 'E27'; Value = '  -2.07%  ' }
    @{ Cell = 'D28'; Value = '''8.56' }
    @{ Cell = 'E28'; Value = '  +0.70%  ' }
    @{ Cell = 'D29'; Value = '''0.130' }
    @{ Cell = 'E29'; Value = '  +3.49%  ' }
    @{ Cell = 'D30'; Value = '''18.29' }
    @{ Cell = 'E30'; Value = '  +1.75%  ' }
    @{ Cell = 'E32'; Value = '  +1.12%  ' }
    @{ Cell = 'D33'; Value = '''0.0567' }
    @{ Cell = 'E33'; Value = '  +0.05%  ' }
    @{ Cell = 'E34'; Value = '  -0.09%  ' }
    @{ Cell = 'B35'; Value = 'InternetComputer(DFINITY)' }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp' }
    @{ Cell = 'D35'; Value = '''4.10' }
    @{ Cell = 'E35'; Value = '  -0.52%  ' }
    @{ Cell = 'B36'; Value = 'ImmutableX' }
    @{ Cell = 'C36'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx' }
    @{ Cell = 'D36'; Value = '''0.918' }
    @{ Cell = 'E36'; Value = '  -5.90%  ' }
    @{ Cell = 'B37'; Value = 'WEMIXToken' }
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix' }
    @{ Cell = 'D37'; Value = '''1.83' }
    @{ Cell = 'E37'; Value = '  +3.56%  ' }
    @{ Cell = 'B38'; Value = 'TrustWalletToken' }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' }
    @{ Cell = 'D38'; Value = '''1.50' }
    @{ Cell = 'E38'; Value = '  +13.37%  ' }
    @{ Cell = 'E39'; Value = '  -1.66%  ' }
    @{ Cell = 'D40'; Value = '''0.0657' }
    @{ Cell = 'E40'; Value = '  +9.87%  ' }
    @{ Cell = 'D41'; Value = '''0.0210' }
    @{ Cell = 'E41'; Value = '  +2.02%  ' }
    @{ Cell = 'E42'; Value = '  -1.86%  ' }
    @{ Cell = 'D43'; Value = '''16.05' }
    @{ Cell = 'E43'; Value = '  +5.42%  ' }
    @{ Cell = 'D44'; Value = '''89.38' }
    @{ Cell = 'E44'; Value = '  -2.05%  ' }
    @{ Cell = 'D45'; Value = '1.351.68' }
    @{ Cell = 'E45'; Value = '  -0.65%  ' }
    @{ Cell = 'E46'; Value = '  +2.17%  ' }
    @{ Cell = 'B47'; Value = 'HuobiToken' }
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht' }
    @{ Cell = 'D47'; Value = '''2.43' }
    @{ Cell = 'E47'; Value = '  +0.17%  ' }
    @{ Cell = 'B48'; Value = 'Gas' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas' }
    @{ Cell = 'D48'; Value = '''12.54' }
    @{ Cell = 'E48'; Value = '  -2.76%  ' }
    @{ Cell = 'E49'; Value = '  +0.39%  ' }
    @{ Cell = 'D50'; Value = '''45.75' }
    @{ Cell = 'E50'; Value = '  -3.59%  ' }
    @{ Cell = 'E51'; Value = '  -2.73%  ' }
)

foreach ($update in $updates) {
    $ws.Range($update.Cell).Value = $update.Value
}
